$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells P1, Q1 -- continue the style of the existing header row
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").Value2 = 14
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value2 = 15

# For rows 2..25, swap column I<->K and M<->O, then add P and Q = 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Range("I$r").Value2
    $kVal = $ws.Range("K$r").Value2
    $mVal = $ws.Range("M$r").Value2
    $oVal = $ws.Range("O$r").Value2

    $ws.Range("I$r").Value2 = $kVal
    $ws.Range("K$r").Value2 = $iVal
    $ws.Range("M$r").Value2 = $oVal
    $ws.Range("O$r").Value2 = $mVal

    $ws.Range("P$r").Value2 = 2
    $ws.Range("Q$r").Value2 = 2
}
